# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# The old scraper only pulled team statistics, not the season record -
# this back-fills Wins/Losses/Ties for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 should look like the rest of the header row
# (bold font, thin border, centered/top aligned) -> copy the format from
# the last existing header cell (AC1) instead of re-building the style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player (rows 2-44) shares the same 2016 KC Royals season record:
# 81 wins, 81 losses, 0 ties.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 81  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 81  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
